# Refactor api response formatting
# Appends a new row (row 94) of data to each of the four worksheets,
# mirroring the existing row 93 layout/formatting.

$wb = $excel.ActiveWorkbook

$rows = @(
    @{
        Sheet = 1
        A = 45880.46178240741
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D = "0x01,0x18"
        E = "0x07"
        F = 400
        G = [double]"5.68631262647113e+23"
        H = 280
        I = 7
    },
    @{
        Sheet = 2
        A = 45880.46178240741
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D = "0x01,0x24"
        E = "0x19"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 292
        I = 25
    },
    @{
        Sheet = 3
        A = 45880.46178240741
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D = "0x00,0x5E"
        E = "0x15"
        F = 110
        G = [double]"5.68631262647113e+23"
        H = 94
        I = 15
    },
    @{
        Sheet = 4
        A = 45880.46178240741
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D = "0x00,0x75"
        E = "0x9"
        F = 130
        G = [double]"5.68631262647113e+23"
        H = 117
        I = 9
    }
)

foreach ($rowInfo in $rows) {
    $ws = $wb.Worksheets.Item($rowInfo.Sheet)
    $newRow = 94
    $srcRow = 93

    # Column A: date/time serial - copy the number format from the row above
    # so the new cell keeps the same style (s="2") as existing date cells.
    $ws.Cells.Item($newRow, 1).Value = $rowInfo.A
    $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($srcRow, 1).NumberFormat

    # Columns B-E: text values (hex byte strings)
    $ws.Cells.Item($newRow, 2).Value = $rowInfo.B
    $ws.Cells.Item($newRow, 3).Value = $rowInfo.C
    $ws.Cells.Item($newRow, 4).Value = $rowInfo.D
    $ws.Cells.Item($newRow, 5).Value = $rowInfo.E

    # Columns F-I: numeric values
    $ws.Cells.Item($newRow, 6).Value = $rowInfo.F
    $ws.Cells.Item($newRow, 7).Value = $rowInfo.G
    $ws.Cells.Item($newRow, 8).Value = $rowInfo.H
    $ws.Cells.Item($newRow, 9).Value = $rowInfo.I
}
